$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '27.092.14'
$ws.Range("E2").Value = '  +0.71%  '
$ws.Range("E3").Value = '  +0.85%  '
$ws.Range("E4").Value = '  +0.06%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '215.92'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.13%  '
$ws.Range("E6").Value = '  -3.45%  '
$ws.Range("E7").Value = '  +0.00%  '
$ws.Range("E8").Value = '  +1.84%  '
$ws.Range("E9").Value = '  +5.49%  '
$ws.Range("E10").Value = '  +0.50%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0890'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -0.65%  '
$ws.Range("D12").Value = '1.919.20'
$ws.Range("E12").Value = '  +0.85%  '
$ws.Range("D13").Value = '1.676.27'
$ws.Range("E13").Value = '  +0.66%  '
$ws.Range("E14").Value = '  +0.94%  '
$ws.Range("E15").Value = '  +2.09%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '66.47'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +0.51%  '
$ws.Range("D17").Value = '27.074.21'
$ws.Range("E17").Value = '  +0.55%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '236.41'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +0.95%  '
$ws.Range("D20").Value = '0.0₃0738'
$ws.Range("E20").Value = '  +0.83%  '
$ws.Range("E21").Value = '  +0.03%  '
$ws.Range("E22").Value = '  +1.38%  '
$ws.Range("E23").Value = '  +1.40%  '
$ws.Range("E24").Value = '  -3.41%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '147.33'
$ws.Range("D25").Style = "Normal"
$ws.Range("E26").Value = '  +2.29%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '16.55'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +4.21%  '
$ws.Range("E28").Value = '  -1.75%  '
$ws.Range("E29").Value = '  +0.19%  '
$ws.Range("E30").Value = '  +0.31%  '
$ws.Range("E31").Value = '  +0.23%  '
$ws.Range("E32").Value = '  +0.56%  '
$ws.Range("D33").Value = '1.547.74'
$ws.Range("E33").Value = '  +6.65%  '
$ws.Range("E34").Value = '  +1.41%  '
$ws.Range("E35").Value = '  +4.72%  '
$ws.Range("E36").Value = '  -1.03%  '
$ws.Range("E37").Value = '  +0.96%  '
$ws.Range("E38").Value = '  +1.34%  '
$ws.Range("E39").Value = '  +2.36%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.05'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +7.34%  '
$ws.Range("E41").Value = '  +0.02%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '67.92'
$ws.Range("D42").Style = "Normal"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '5.54'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -3.47%  '
$ws.Range("E44").Value = '  -0.99%  '
$ws.Range("D45").Value = '1.823.17'
$ws.Range("E45").Value = '  +0.76%  '
$ws.Range("E46").Value = '  -0.37%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '90.60'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +0.11%  '
$ws.Range("E48").Value = '  +3.24%  '
$ws.Range("E49").Value = '  +0.66%  '
$ws.Range("E50").Value = '  +1.93%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '8.01'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +6.53%  '
